# Auto-generated: update scheduled market-price derived cells across sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1908.5
$ws.Range("I2").Value = 1518.3334
$ws.Range("J2").Value = 2298.6667
$ws.Range("K2").Value = 1518.3334
$ws.Range("L2").Value = 2298.6667
$ws.Range("M2").Value = -1405.3334
$ws.Range("N2").Value = -2524.6667
$ws.Range("H69").Value = 7999.6665
$ws.Range("J69").Value = 9999.5
$ws.Range("L69").Value = 29998.5
$ws.Range("N69").Value = -31746.5
$ws.Range("H72").Value = 7999.6665
$ws.Range("J72").Value = 9999.5
$ws.Range("L72").Value = 89995.5
$ws.Range("N72").Value = -98731.5
$ws.Range("H80").Value = 9304.764999999999
$ws.Range("J80").Value = 13749.909
$ws.Range("L80").Value = 41249.727
$ws.Range("N80").Value = -43245.727
$ws.Range("H83").Value = 9304.764999999999
$ws.Range("J83").Value = 13749.909
$ws.Range("L83").Value = 123749.181
$ws.Range("N83").Value = -133733.181
$ws.Range("H101").Value = 197863.75
$ws.Range("J101").Value = 394992.5
$ws.Range("L101").Value = 1184977.5
$ws.Range("N101").Value = -1188221.5
$ws.Range("H125").Value = 937.9375
$ws.Range("J125").Value = 968.9167
$ws.Range("L125").Value = 8720.2503
$ws.Range("N125").Value = -13640.2503
$ws.Range("H132").Value = 1570.9678
$ws.Range("I132").Value = 1123.1333
$ws.Range("K132").Value = 3369.3999
$ws.Range("M132").Value = -839.3998999999999
$ws.Range("H133").Value = 93574
$ws.Range("J133").Value = 93574
$ws.Range("L133").Value = 93574
$ws.Range("N133").Value = -103694
$ws.Range("H134").Value = 80525.11
$ws.Range("J134").Value = 91135.86
$ws.Range("L134").Value = 91135.86
$ws.Range("N134").Value = -101275.86
$ws.Range("H136").Value = 72869.8
$ws.Range("J136").Value = 82087.25
$ws.Range("L136").Value = 82087.25
$ws.Range("N136").Value = -92287.25
$ws.Range("H138").Value = 2593.5825
$ws.Range("I138").Value = 1519.0312
$ws.Range("J138").Value = 3176.39
$ws.Range("K138").Value = 4557.0936
$ws.Range("L138").Value = 9529.17
$ws.Range("M138").Value = 582.9063999999998
$ws.Range("N138").Value = -19809.17
$ws.Range("H139").Value = 99990
$ws.Range("J139").Value = 99990
$ws.Range("L139").Value = 99990
$ws.Range("N139").Value = -110270
$ws.Range("H140").Value = 91991.664
$ws.Range("J140").Value = 91991.664
$ws.Range("L140").Value = 91991.664
$ws.Range("N140").Value = -102351.664

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 7230.625
$ws.Range("I26").Value = 5405.2856
$ws.Range("K26").Value = 5405.2856
$ws.Range("M26").Value = -5075.2856
$ws.Range("H32").Value = 7838.808
$ws.Range("I32").Value = 3823.0625
$ws.Range("K32").Value = 3823.0625
$ws.Range("M32").Value = -3536.0625
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("H36").Value = 14500
$ws.Range("I36").Value = 24000
$ws.Range("K36").Value = 24000
$ws.Range("M36").Value = -23654
$ws.Range("H74").Value = 44558.668
$ws.Range("I74").Value = 47017.273
$ws.Range("K74").Value = 47017.273
$ws.Range("M74").Value = -46143.273
$ws.Range("H77").Value = 44558.668
$ws.Range("I77").Value = 47017.273
$ws.Range("K77").Value = 235086.365
$ws.Range("M77").Value = -230718.365
$ws.Range("H88").Value = 1668980.9
$ws.Range("J88").Value = 2945
$ws.Range("L88").Value = 2945
$ws.Range("N88").Value = -3757
$ws.Range("H91").Value = 1668980.9
$ws.Range("J91").Value = 2945
$ws.Range("L91").Value = 2945
$ws.Range("N91").Value = -5753
$ws.Range("H122").Value = 3837.0217
$ws.Range("I122").Value = 3823.7837
$ws.Range("J122").Value = 3891.4443
$ws.Range("K122").Value = 11471.3511
$ws.Range("L122").Value = 11674.3329
$ws.Range("M122").Value = -9021.3511
$ws.Range("N122").Value = -16574.3329

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2520.6924
$ws.Range("I105").Value = 1773.4546
$ws.Range("K105").Value = 1773.4546
$ws.Range("M105").Value = -26.45460000000003
$ws.Range("H117").Value = 85358
$ws.Range("J117").Value = 85358
$ws.Range("L117").Value = 85358
$ws.Range("N117").Value = -94536
$ws.Range("H134").Value = 1593.2759
$ws.Range("I134").Value = 1158.5834
$ws.Range("K134").Value = 3475.7502
$ws.Range("M134").Value = -940.7501999999999
$ws.Range("H138").Value = 73705.57000000001
$ws.Range("J138").Value = 73705.57000000001
$ws.Range("L138").Value = 73705.57000000001
$ws.Range("N138").Value = -83985.57000000001
$ws.Range("H140").Value = 56710.816
$ws.Range("J140").Value = 56710.816
$ws.Range("L140").Value = 56710.816
$ws.Range("N140").Value = -67070.81599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4020.8474
$ws.Range("I31").Value = 2075.375
$ws.Range("K31").Value = 2075.375
$ws.Range("M31").Value = -1780.375
$ws.Range("H34").Value = 4020.8474
$ws.Range("I34").Value = 2075.375
$ws.Range("K34").Value = 2075.375
$ws.Range("M34").Value = -1873.375
$ws.Range("H132").Value = 2957
$ws.Range("I132").Value = 1933.3334
$ws.Range("K132").Value = 5800.0002
$ws.Range("M132").Value = -3270.0002
$ws.Range("H134").Value = 35720.965
$ws.Range("I134").Value = 2363.238
$ws.Range("K134").Value = 7089.714
$ws.Range("M134").Value = -4554.714
$ws.Range("H138").Value = 55282.5
$ws.Range("J138").Value = 54894.285
$ws.Range("L138").Value = 54894.285
$ws.Range("N138").Value = -65174.285

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 163.6
$ws.Range("I2").Value = 14.5
$ws.Range("J2").Value = 263
$ws.Range("K2").Value = 87
$ws.Range("L2").Value = 1578
$ws.Range("M2").Value = 26
$ws.Range("N2").Value = -1804
$ws.Range("H55").Value = 111122110
$ws.Range("J55").Value = 125012250
$ws.Range("L55").Value = 375036750
$ws.Range("N55").Value = -375037104

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1860.25
$ws.Range("I113").Value = 1353.1666
$ws.Range("K113").Value = 1353.1666
$ws.Range("M113").Value = 816.8334
$ws.Range("H122").Value = 12514.25
$ws.Range("I122").Value = 3611.5715
$ws.Range("K122").Value = 10834.7145
$ws.Range("M122").Value = -8384.7145
$ws.Range("H132").Value = 2117.3513
$ws.Range("I132").Value = 1924.3
$ws.Range("J132").Value = 2944.7144
$ws.Range("K132").Value = 5772.9
$ws.Range("L132").Value = 8834.143199999999
$ws.Range("M132").Value = -3242.9
$ws.Range("N132").Value = -13894.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2648566
$ws.Range("I40").Value = 2877.8333
$ws.Range("J40").Value = 6176150.5
$ws.Range("K40").Value = 2877.8333
$ws.Range("L40").Value = 6176150.5
$ws.Range("M40").Value = -2741.8333
$ws.Range("N40").Value = -6176422.5
$ws.Range("H132").Value = 2809.2593
$ws.Range("I132").Value = 2368
$ws.Range("K132").Value = 7104
$ws.Range("M132").Value = -4574
$ws.Range("H136").Value = 5692.346
$ws.Range("J136").Value = 5651.846
$ws.Range("L136").Value = 16955.538
$ws.Range("N136").Value = -22055.538

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2357.5715
$ws.Range("J81").Value = 2450.3333
$ws.Range("L81").Value = 4900.6666
$ws.Range("N81").Value = -7022.6666
$ws.Range("H84").Value = 2357.5715
$ws.Range("J84").Value = 2450.3333
$ws.Range("L84").Value = 24503.333
$ws.Range("N84").Value = -35111.333
$ws.Range("H113").Value = 1826.1111
$ws.Range("I113").Value = 1916.3334
$ws.Range("K113").Value = 5749.0002
$ws.Range("M113").Value = -3579.0002
$ws.Range("H122").Value = 4444.55
$ws.Range("I122").Value = 3036.375
$ws.Range("J122").Value = 5383.3335
$ws.Range("K122").Value = 9109.125
$ws.Range("L122").Value = 16150.0005
$ws.Range("M122").Value = -6659.125
$ws.Range("N122").Value = -21050.0005
